$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 26823
$ws.Cells.Item(4, 6).Value = 591
$ws.Cells.Item(5, 6).Value = 258
$ws.Cells.Item(7, 6).Value = 178
$ws.Cells.Item(11, 6).Value = 449
$ws.Cells.Item(12, 6).Value = 191
$ws.Cells.Item(15, 6).Value = 74
$ws.Cells.Item(16, 6).Value = 440
$ws.Cells.Item(18, 6).Value = 1564
$ws.Cells.Item(20, 6).Value = 54
$ws.Cells.Item(21, 6).Value = 444
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 4512
$ws.Cells.Item(3, 6).Value = 237
$ws.Cells.Item(6, 6).Value = 204
$ws.Cells.Item(7, 6).Value = 204
$ws.Cells.Item(8, 6).Value = 37
$ws.Cells.Item(15, 6).Value = 64
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 5104
$ws.Cells.Item(3, 6).Value = 244
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 5104
$ws.Cells.Item(4, 6).Value = 244
$ws.Cells.Item(5, 6).Value = 26823
$ws.Cells.Item(6, 6).Value = 591
$ws.Cells.Item(7, 6).Value = 4512
$ws.Cells.Item(8, 6).Value = 258
$ws.Cells.Item(9, 6).Value = 237
$ws.Cells.Item(13, 6).Value = 178
$ws.Cells.Item(14, 6).Value = 204
$ws.Cells.Item(15, 6).Value = 204
$ws.Cells.Item(16, 6).Value = 37
$ws.Cells.Item(23, 6).Value = 449
$ws.Cells.Item(24, 6).Value = 191
$ws.Cells.Item(28, 6).Value = 74
$ws.Cells.Item(31, 6).Value = 440
$ws.Cells.Item(33, 6).Value = 64
$ws.Cells.Item(34, 6).Value = 1564
$ws.Cells.Item(37, 6).Value = 54
$ws.Cells.Item(38, 6).Value = 444
